# Update "想去人数" (want-to-go count) figures in the F column of the
# "展览" and "全部类型" sheets to match the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1027
    "F3"  = 298
    "F4"  = 1419
    "F5"  = 8543
    "F7"  = 481
    "F11" = 3440
    "F14" = 64
    "F15" = 997
    "F17" = 1095
    "F18" = 297
    "F19" = 170
    "F20" = 2123
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
